$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1609.25
$ws.Range("I38").Value = 959.5454999999999
$ws.Range("J38").Value = 3038.6
$ws.Range("K38").Value = 2878.6365
$ws.Range("L38").Value = 9115.799999999999
$ws.Range("M38").Value = -2506.6365
$ws.Range("N38").Value = -9859.799999999999

$ws.Range("H129").Value = 4063.3333
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 1963.75
$ws.Range("I132").Value = 1887.7916
$ws.Range("J132").Value = 2419.5
$ws.Range("K132").Value = 5663.3748
$ws.Range("L132").Value = 7258.5
$ws.Range("M132").Value = -3133.3748
$ws.Range("N132").Value = -12318.5

$ws.Range("H137").Value = 1685.2
$ws.Range("I137").Value = 1718.7778
$ws.Range("K137").Value = 5156.3334
$ws.Range("M137").Value = -2606.3334

$ws.Range("H138").Value = 5435.25
$ws.Range("J138").Value = 4202.3
$ws.Range("L138").Value = 12606.9
$ws.Range("N138").Value = -22886.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12804.8
$ws.Range("J2").Value = 2680
$ws.Range("L2").Value = 2680
$ws.Range("N2").Value = -2906

$ws.Range("H4").Value = 450.6
$ws.Range("I4").Value = 450.6
$ws.Range("K4").Value = 450.6
$ws.Range("M4").Value = -334.6

$ws.Range("H102").Value = 1561.25
$ws.Range("I102").Value = 1603.1818
$ws.Range("K102").Value = 1603.1818
$ws.Range("M102").Value = 18.81819999999993

$ws.Range("H110").Value = 3358.1904
$ws.Range("I110").Value = 2758.4285
$ws.Range("J110").Value = 4557.7144
$ws.Range("K110").Value = 2758.4285
$ws.Range("L110").Value = 4557.7144
$ws.Range("M110").Value = -713.4285
$ws.Range("N110").Value = -8647.714400000001

$ws.Range("H116").Value = 12804.8
$ws.Range("J116").Value = 2680
$ws.Range("L116").Value = 2680
$ws.Range("N116").Value = -7268

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12804.8
$ws.Range("J3").Value = 2680
$ws.Range("L3").Value = 2680
$ws.Range("N3").Value = -2908

$ws.Range("H22").Value = 325.375
$ws.Range("I22").Value = 210
$ws.Range("J22").Value = 363.83334
$ws.Range("K22").Value = 210
$ws.Range("L22").Value = 363.83334
$ws.Range("M22").Value = -37
$ws.Range("N22").Value = -709.83334

$ws.Range("H105").Value = 2164.611
$ws.Range("I105").Value = 1664.2667
$ws.Range("K105").Value = 1664.2667
$ws.Range("M105").Value = 82.7333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1564.4517
$ws.Range("I31").Value = 1990.7273
$ws.Range("J31").Value = 1330
$ws.Range("K31").Value = 1990.7273
$ws.Range("L31").Value = 1330
$ws.Range("M31").Value = -1695.7273
$ws.Range("N31").Value = -1920

$ws.Range("H33").Value = 7358.7
$ws.Range("I33").Value = 3365.1428
$ws.Range("J33").Value = 16677
$ws.Range("K33").Value = 3365.1428
$ws.Range("L33").Value = 16677
$ws.Range("M33").Value = -2986.1428
$ws.Range("N33").Value = -17435

$ws.Range("H34").Value = 1564.4517
$ws.Range("I34").Value = 1990.7273
$ws.Range("J34").Value = 1330
$ws.Range("K34").Value = 1990.7273
$ws.Range("L34").Value = 1330
$ws.Range("M34").Value = -1788.7273
$ws.Range("N34").Value = -1734

$ws.Range("H59").Value = 28476.6
$ws.Range("J59").Value = 26224
$ws.Range("L59").Value = 26224
$ws.Range("N59").Value = -28514

$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -39492

$ws.Range("H107").Value = 2572
$ws.Range("I107").Value = 2141.8333
$ws.Range("K107").Value = 2141.8333
$ws.Range("M107").Value = -221.8332999999998

$ws.Range("H122").Value = 2003.9474
$ws.Range("I122").Value = 1974.8572
$ws.Range("K122").Value = 5924.571599999999
$ws.Range("M122").Value = -3474.571599999999

$ws.Range("H132").Value = 1450.9565
$ws.Range("I132").Value = 1450.9565
$ws.Range("K132").Value = 4352.8695
$ws.Range("M132").Value = -1822.8695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 127568
$ws.Range("I4").Value = 524.5
$ws.Range("J4").Value = 169915.83
$ws.Range("K4").Value = 1573.5
$ws.Range("L4").Value = 509747.49
$ws.Range("M4").Value = -1461.5
$ws.Range("N4").Value = -509971.49

$ws.Range("H11").Value = 444.66666
$ws.Range("I11").Value = 445
$ws.Range("J11").Value = 444.33334
$ws.Range("K11").Value = 1335
$ws.Range("L11").Value = 1333.00002
$ws.Range("M11").Value = -1195
$ws.Range("N11").Value = -1613.00002

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H36").Value = 401
$ws.Range("I36").Value = 401.25
$ws.Range("J36").Value = 400
$ws.Range("K36").Value = 1203.75
$ws.Range("L36").Value = 1200
$ws.Range("M36").Value = -1034.75
$ws.Range("N36").Value = -1538

$ws.Range("H47").Value = 1476.75
$ws.Range("I47").Value = 290.125
$ws.Range("J47").Value = 3850
$ws.Range("K47").Value = 870.375
$ws.Range("L47").Value = 11550
$ws.Range("M47").Value = -439.375
$ws.Range("N47").Value = -12412

$ws.Range("H61").Value = 2661.5
$ws.Range("I61").Value = 2392.8
$ws.Range("J61").Value = 4005
$ws.Range("K61").Value = 7178.400000000001
$ws.Range("L61").Value = 12015
$ws.Range("M61").Value = -6963.400000000001
$ws.Range("N61").Value = -12445

$ws.Range("H122").Value = 334.7143
$ws.Range("J122").Value = 300.72726
$ws.Range("L122").Value = 2706.54534
$ws.Range("N122").Value = -7606.545340000001

$ws.Range("H132").Value = 4433.8125
$ws.Range("J132").Value = 4912
$ws.Range("L132").Value = 44208
$ws.Range("N132").Value = -49268

$ws.Range("H137").Value = 1500505.8
$ws.Range("I137").Value = 1264.5
$ws.Range("K137").Value = 3793.5
$ws.Range("M137").Value = 1306.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 19291
$ws.Range("I35").Value = 14632.5
$ws.Range("J35").Value = 23949.5
$ws.Range("K35").Value = 14632.5
$ws.Range("L35").Value = 23949.5
$ws.Range("M35").Value = -14334.5
$ws.Range("N35").Value = -24545.5

$ws.Range("H113").Value = 2057.5386
$ws.Range("I113").Value = 2165
$ws.Range("K113").Value = 2165
$ws.Range("M113").Value = 5

$ws.Range("H122").Value = 2258
$ws.Range("I122").Value = 2221
$ws.Range("J122").Value = 2332
$ws.Range("K122").Value = 6663
$ws.Range("L122").Value = 6996
$ws.Range("M122").Value = -4213
$ws.Range("N122").Value = -11896

$ws.Range("H132").Value = 2310.0527
$ws.Range("I132").Value = 2310.0527
$ws.Range("K132").Value = 6930.158100000001
$ws.Range("M132").Value = -4400.158100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 254025
$ws.Range("J20").Value = 254025
$ws.Range("L20").Value = 254025
$ws.Range("N20").Value = -254477

$ws.Range("H40").Value = 3293.5
$ws.Range("I40").Value = 2709.5
$ws.Range("K40").Value = 2709.5
$ws.Range("M40").Value = -2573.5

$ws.Range("H46").Value = 1935.6364
$ws.Range("I46").Value = 2141.7144
$ws.Range("J46").Value = 1575
$ws.Range("K46").Value = 2141.7144
$ws.Range("L46").Value = 1575
$ws.Range("M46").Value = -1953.7144
$ws.Range("N46").Value = -1951

$ws.Range("H61").Value = 1579.5
$ws.Range("J61").Value = 1442.5
$ws.Range("L61").Value = 1442.5
$ws.Range("N61").Value = -1846.5

$ws.Range("H113").Value = 1579.5
$ws.Range("J113").Value = 1442.5
$ws.Range("L113").Value = 1442.5
$ws.Range("N113").Value = -5782.5

$ws.Range("H122").Value = 3367.6667
$ws.Range("I122").Value = 3367.6667
$ws.Range("K122").Value = 10103.0001
$ws.Range("M122").Value = -7653.000100000001

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 2510.1
$ws.Range("I132").Value = 2421.7778
$ws.Range("K132").Value = 7265.3334
$ws.Range("M132").Value = -4735.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1330.6666
$ws.Range("I113").Value = 1330.6666
$ws.Range("K113").Value = 3991.9998
$ws.Range("M113").Value = -1821.9998

$ws.Range("H126").Value = 3235.0527
$ws.Range("I126").Value = 2433.4285
$ws.Range("J126").Value = 5479.6
$ws.Range("K126").Value = 7300.2855
$ws.Range("L126").Value = 16438.8
$ws.Range("M126").Value = -4830.2855
$ws.Range("N126").Value = -21378.8

$ws.Range("H132").Value = 2541.5557
$ws.Range("I132").Value = 2541.5557
$ws.Range("K132").Value = 7624.6671
$ws.Range("M132").Value = -5094.6671
